$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "23.834.60"
Set-TextValue $ws.Range("E2") "  -3.07%  "
Set-TextValue $ws.Range("D3") "1.615.51"
Set-TextValue $ws.Range("E3") "  -3.37%  "
Set-TextValue $ws.Range("E4") "  -0.07%  "
Set-TextValue $ws.Range("D5") "307.70"
Set-TextValue $ws.Range("E5") "  -2.12%  "
Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("E6") "  -0.07%  "
Set-TextValue $ws.Range("D7") "0.3924"
Set-TextValue $ws.Range("E7") "  -0.86%  "
Set-TextValue $ws.Range("D8") "0.3831"
Set-TextValue $ws.Range("E8") "  -3.24%  "
Set-TextValue $ws.Range("D9") "1.001"
Set-TextValue $ws.Range("E9") "  -0.12%  "
Set-TextValue $ws.Range("D10") "1.353"
Set-TextValue $ws.Range("E10") "  -3.14%  "
Set-TextValue $ws.Range("D11") "49.33"
Set-TextValue $ws.Range("E11") "  -2.08%  "
Set-TextValue $ws.Range("D12") "0.08427"
Set-TextValue $ws.Range("E12") "  -2.53%  "
Set-TextValue $ws.Range("D13") "23.59"
Set-TextValue $ws.Range("E13") "  -7.32%  "
Set-TextValue $ws.Range("D14") "7.033"
Set-TextValue $ws.Range("E14") "  -3.75%  "
Set-TextValue $ws.Range("D15") "7.548"
Set-TextValue $ws.Range("E15") "  -1.85%  "
Set-TextValue $ws.Range("D16") "0.00001278"
Set-TextValue $ws.Range("D17") "1.618.29"
Set-TextValue $ws.Range("E17") "  -3.61%  "
Set-TextValue $ws.Range("D18") "93.70"
Set-TextValue $ws.Range("E18") "  -0.22%  "
Set-TextValue $ws.Range("D19") "0.06925"
Set-TextValue $ws.Range("E19") "  -1.33%  "
Set-TextValue $ws.Range("D20") "19.98"
Set-TextValue $ws.Range("E20") "  -5.87%  "
Set-TextValue $ws.Range("D21") "6.796"
Set-TextValue $ws.Range("E21") "  -4.04%  "
Set-TextValue $ws.Range("E22") "  -0.05%  "
Set-TextValue $ws.Range("D23") "13.40"
Set-TextValue $ws.Range("E23") "  -3.64%  "
Set-TextValue $ws.Range("D24") "23.847.02"
Set-TextValue $ws.Range("E24") "  -3.03%  "
Set-TextValue $ws.Range("E25") "  +4.64%  "
Set-TextValue $ws.Range("D26") "2.837"
Set-TextValue $ws.Range("E26") "  +2.46%  "
Set-TextValue $ws.Range("D27") "22.15"
Set-TextValue $ws.Range("E27") "  -3.74%  "
Set-TextValue $ws.Range("D28") "156.74"
Set-TextValue $ws.Range("E28") "  -1.99%  "
Set-TextValue $ws.Range("D29") "139.49"
Set-TextValue $ws.Range("E29") "  -4.21%  "
Set-TextValue $ws.Range("E30") "  -9.61%  "
Set-TextValue $ws.Range("D31") "7.770"
Set-TextValue $ws.Range("E31") "  -6.45%  "
Set-TextValue $ws.Range("D32") "2.486"
Set-TextValue $ws.Range("E32") "  -1.74%  "
Set-TextValue $ws.Range("D33") "1.795.20"
Set-TextValue $ws.Range("E33") "  -3.93%  "
Set-TextValue $ws.Range("D34") "0.08069"
Set-TextValue $ws.Range("E34") "  -2.24%  "
Set-TextValue $ws.Range("D35") "0.9719"
Set-TextValue $ws.Range("E35") "  -2.08%  "
Set-TextValue $ws.Range("D36") "0.02874"
Set-TextValue $ws.Range("E36") "  -6.76%  "
Set-TextValue $ws.Range("D37") "6.546"
Set-TextValue $ws.Range("E37") "  -5.35%  "
Set-TextValue $ws.Range("D38") "0.2659"
Set-TextValue $ws.Range("E38") "  -5.31%  "
Set-TextValue $ws.Range("D39") "0.09109"
Set-TextValue $ws.Range("E39") "  -5.58%  "
Set-TextValue $ws.Range("D40") "10.30"
Set-TextValue $ws.Range("E40") "  -0.16%  "
Set-TextValue $ws.Range("D41") "13.57"
Set-TextValue $ws.Range("E41") "  +0.15%  "
Set-TextValue $ws.Range("D42") "1.424"
Set-TextValue $ws.Range("E42") "  -6.38%  "
Set-TextValue $ws.Range("D43") "0.7475"
Set-TextValue $ws.Range("E43") "  -5.09%  "
Set-TextValue $ws.Range("D44") "16.05"
Set-TextValue $ws.Range("E44") "  -3.17%  "
Set-TextValue $ws.Range("D45") "0.6894"
Set-TextValue $ws.Range("E45") "  -2.82%  "
Set-TextValue $ws.Range("D46") "2.460"
Set-TextValue $ws.Range("E46") "  -4.03%  "
Set-TextValue $ws.Range("D47") "4.068"
Set-TextValue $ws.Range("E47") "  -2.72%  "
Set-TextValue $ws.Range("D48") "1.001"
Set-TextValue $ws.Range("E48") "  -0.13%  "
Set-TextValue $ws.Range("D49") "0.08216"
Set-TextValue $ws.Range("E49") "  -4.84%  "
Set-TextValue $ws.Range("D50") "134.53"
Set-TextValue $ws.Range("E50") "  -2.52%  "
Set-TextValue $ws.Range("D51") "1.203"
Set-TextValue $ws.Range("E51") "  -9.35%  "
